$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: first affectee record ---------------------------------------
$ws.Range("A2").Value = "Name1"
$ws.Range("B2").Value = "Father Name 1"
$ws.Range("C2").Value = 12234123
$ws.Range("D2").Value = "Address of affectee"
$ws.Range("E2").Value = "Peshawar"
$ws.Range("F2").Value = "Some Reason"
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = "big"
$ws.Range("J2").Value = "2015-31-04"
$ws.Range("K2").Value = "yes"
$ws.Range("L2").Value = "yes"
$ws.Range("M2").Value = "yes"
$ws.Range("N2").Value = "yes"
$ws.Range("O2").Value = "yes"
$ws.Range("P2").Value = "yes"

# --- Row 3: second affectee record ---------------------------------------
$ws.Range("A3").Value = "Name2"
$ws.Range("B3").Value = "Father Name 2"
$ws.Range("C3").Value = 12234123
$ws.Range("D3").Value = "Address of affectee"
$ws.Range("E3").Value = "Peshawar"
$ws.Range("F3").Value = "Some Reason"
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = "big"
$ws.Range("J3").Value = "2015-31-04"
$ws.Range("K3").Value = "yes"
$ws.Range("L3").Value = "yes"
$ws.Range("M3").Value = "yes"
$ws.Range("N3").Value = "yes"
$ws.Range("O3").Value = "yes"
$ws.Range("P3").Value = "yes"

# --- Column I ("Date Of Incident") formatted as Text, filled in last -----
$ws.Columns.Item(9).NumberFormat = "@"
$ws.Range("I1").NumberFormat = "@"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2012-03-02"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2012-03-02"

# --- Column C width --------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 18.5

# --- View: scroll back to top-left, move selection to I7 -------------------
$ws.Range("I7").Select()

# --- Page setup: portrait orientation --------------------------------------
$ws.PageSetup.Orientation = 1
